# Helper: write a value into a cell while forcing it to be stored as TEXT
# (even if it looks numeric, e.g. "17.06" or "009011"), without leaving a
# lingering custom number-format/style on the cell.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# The existing "总计" sheet (5th tab) will become "2022-Q1": first make an
# exact copy of it (values + formatting) right after itself -- that copy
# will keep being the "总计" sheet, while the original tab gets repurposed.
$totalOriginal = $wb.Worksheets.Item(5)
$totalOriginal.Copy($null, $totalOriginal)

$newTotal = $wb.Worksheets.Item(6)
$newTotal.Name = "总计-new"
$totalOriginal.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 1) Rebuild the "2022-Q1" sheet (formerly "总计") with the fund holdings
#    table (same shape as the other quarter sheets).
# ---------------------------------------------------------------------
$ws1 = $totalOriginal

# Clear everything first
$ws1.Cells.Clear()

# Header row
$ws1.Cells.Item(1,2).Value = "基金代码"
$ws1.Cells.Item(1,3).Value = "基金名称"
$ws1.Cells.Item(1,4).Value = "基金规模"
$ws1.Cells.Item(1,5).Value = "股票总仓位"
$ws1.Cells.Item(1,6).Value = "仓位占比"
$ws1.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws1.Cells.Item(1,8).Value = "仓位排名"
$ws1.Range("B1:H1").Font.Bold = $true
$ws1.Range("B1:H1").Borders.LineStyle = 1
$ws1.Range("B1:H1").HorizontalAlignment = -4108
$ws1.Range("B1:H1").VerticalAlignment = -4160

# Row 2
$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,1).Font.Bold = $true
$ws1.Cells.Item(2,1).Borders.LineStyle = 1
$ws1.Cells.Item(2,1).HorizontalAlignment = -4108
$ws1.Cells.Item(2,1).VerticalAlignment = -4160
Set-TextValue $ws1.Cells.Item(2,2) "009011"
Set-TextValue $ws1.Cells.Item(2,3) "华夏睿阳一年持有期混合"
Set-TextValue $ws1.Cells.Item(2,4) "17.06"
Set-TextValue $ws1.Cells.Item(2,5) "82.70"
Set-TextValue $ws1.Cells.Item(2,6) "2.10"
Set-TextValue $ws1.Cells.Item(2,7) "0.3583"
$ws1.Cells.Item(2,8).Value = 7

# Row 3
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,1).Font.Bold = $true
$ws1.Cells.Item(3,1).Borders.LineStyle = 1
$ws1.Cells.Item(3,1).HorizontalAlignment = -4108
$ws1.Cells.Item(3,1).VerticalAlignment = -4160
Set-TextValue $ws1.Cells.Item(3,2) "004332"
Set-TextValue $ws1.Cells.Item(3,3) "恒生前海沪港深新兴产业精选混合"
Set-TextValue $ws1.Cells.Item(3,4) "0.52"
Set-TextValue $ws1.Cells.Item(3,5) "80.98"
Set-TextValue $ws1.Cells.Item(3,6) "4.53"
Set-TextValue $ws1.Cells.Item(3,7) "0.0236"
$ws1.Cells.Item(3,8).Value = 3

# Row 4
$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,1).Font.Bold = $true
$ws1.Cells.Item(4,1).Borders.LineStyle = 1
$ws1.Cells.Item(4,1).HorizontalAlignment = -4108
$ws1.Cells.Item(4,1).VerticalAlignment = -4160
Set-TextValue $ws1.Cells.Item(4,2) "011800"
Set-TextValue $ws1.Cells.Item(4,3) "申万菱信价值精选混合型证券投资基金"
Set-TextValue $ws1.Cells.Item(4,4) "0.57"
Set-TextValue $ws1.Cells.Item(4,5) "81.46"
Set-TextValue $ws1.Cells.Item(4,6) "3.60"
Set-TextValue $ws1.Cells.Item(4,7) "0.0205"
$ws1.Cells.Item(4,8).Value = 3

# ---------------------------------------------------------------------
# 2) Rebuild the new "总计" sheet: same as the old one, plus a 2022-Q1
#    row inserted right after the header.
# ---------------------------------------------------------------------
$ws2 = $newTotal

# Shift the 4 existing data rows down by one (from the bottom up)
for ($r = 5; $r -ge 2; $r--) {
    $destRow = $r + 1
    $ws2.Range("A$r:D$r").Copy()
    $ws2.Range("A${destRow}:D${destRow}").PasteSpecial(-4122)
    $ws2.Cells.Item($destRow, 1).Value = $r - 1
    $ws2.Cells.Item($destRow, 2).Value = $ws2.Cells.Item($r, 2).Value()
    $ws2.Cells.Item($destRow, 3).Value = $ws2.Cells.Item($r, 3).Value()
    $ws2.Cells.Item($destRow, 4).Value = $ws2.Cells.Item($r, 4).Value()
}

# New row 2: 2022-Q1 summary
$ws2.Range("A2:D2").Copy()
$ws2.Range("A3:D3").PasteSpecial(-4122)
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "2022-Q1"
$ws2.Cells.Item(2,3).Value = 3
$ws2.Cells.Item(2,4).Value = 0.4

$ws2.Name = "总计"

Write-Output "done"
